$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.366.48"
$ws.Range("E2").Value = "  -0.39%  "

$ws.Range("D3").Value = "2.982.31"
$ws.Range("E3").Value = "  +1.20%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.22%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "383.58"
$ws.Range("E5").Value = "  +2.38%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "103.41"
$ws.Range("E6").Value = "  -1.47%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.544"
$ws.Range("E7").Value = "  -0.82%  "

$ws.Range("E8").Value = "  +0.14%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.591"
$ws.Range("E9").Value = "  -0.85%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "37.02"
$ws.Range("E10").Value = "  -1.15%  "

$ws.Range("E11").Value = "  -0.24%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0844"
$ws.Range("E12").Value = "  +0.13%  "

$ws.Range("D13").Value = "3.456.20"
$ws.Range("E13").Value = "  +1.81%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "18.25"
$ws.Range("E14").Value = "  -1.34%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.53"
$ws.Range("E15").Value = "  +1.26%  "

$ws.Range("D16").Value = "2.982.98"
$ws.Range("E16").Value = "  +1.77%  "

$ws.Range("E17").Value = "  +6.70%  "

$ws.Range("D18").Value = "51.368.69"
$ws.Range("E18").Value = "  -0.25%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.26"
$ws.Range("E19").Value = "  -1.44%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.40"
$ws.Range("E20").Value = "  +0.63%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.84"
$ws.Range("E21").Value = "  -2.20%  "

$ws.Range("D22").Value = "0.0₃0958"
$ws.Range("E22").Value = "  +0.64%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "69.05"
$ws.Range("E23").Value = "  +0.29%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "262.04"
$ws.Range("E24").Value = "  -0.08%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.91"
$ws.Range("E25").Value = "  +7.30%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.19"
$ws.Range("E26").Value = "  +13.39%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.55"
$ws.Range("E27").Value = "  +10.65%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.118"
$ws.Range("E28").Value = "  +13.78%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.169"
$ws.Range("E29").Value = "  -2.00%  "

$ws.Range("E30").Value = "  -0.20%  "

$ws.Range("E31").Value = "  +0.02%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "25.92"
$ws.Range("E32").Value = "  -0.27%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "9.86"
$ws.Range("E33").Value = "  -1.04%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "34.61"
$ws.Range("E34").Value = "  -0.98%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "50.90"
$ws.Range("E35").Value = "  -0.41%  "

$ws.Range("B36").Value = "Toncoin"
$ws.Range("C36").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.07"
$ws.Range("E36").Value = "  -2.25%  "

$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0455"
$ws.Range("E37").Value = "  +6.03%  "

$ws.Range("E38").Value = "  -0.04%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.99"
$ws.Range("E39").Value = "  -1.44%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "16.98"
$ws.Range("E40").Value = "  -1.73%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.57"
$ws.Range("E41").Value = "  -0.85%  "

$ws.Range("E42").Value = "  +1.09%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.82"
$ws.Range("E43").Value = "  -2.24%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "122.79"
$ws.Range("E44").Value = "  +2.57%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "21.54"
$ws.Range("E45").Value = "  -2.74%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.09"
$ws.Range("E46").Value = "  -0.69%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.273"
$ws.Range("E47").Value = "  +4.85%  "

$ws.Range("E48").Value = "  +2.52%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.28"
$ws.Range("E49").Value = "  +1.81%  "

$ws.Range("D50").Value = "2.032.23"
$ws.Range("E50").Value = "  -0.18%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0332"
$ws.Range("E51").Value = "  +2.30%  "
